$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values that look numeric are kept as text (matching original inline-string formatting)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '40.333.25'
$ws.Range("E2").Value = '  +0.69%  '
$ws.Range("D3").Value = '2.220.34'
$ws.Range("E3").Value = '  -0.29%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '299.74'
$ws.Range("E5").Value = '  +1.83%  '
$ws.Range("D6").Value = '88.39'
$ws.Range("E6").Value = '  +1.24%  '
$ws.Range("E7").Value = '  +0.80%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D10").Value = '52.69'
$ws.Range("E10").Value = '  +7.73%  '
$ws.Range("D11").Value = '31.15'
$ws.Range("E11").Value = '  +3.24%  '
$ws.Range("D12").Value = '0.0786'
$ws.Range("E12").Value = '  +0.82%  '
$ws.Range("E13").Value = '  +2.53%  '
$ws.Range("D14").Value = '6.43'
$ws.Range("E14").Value = '  -0.69%  '
$ws.Range("D15").Value = '2.561.20'
$ws.Range("E15").Value = '  -0.26%  '
$ws.Range("D16").Value = '13.90'
$ws.Range("E16").Value = '  +0.89%  '
$ws.Range("D17").Value = '2.204.26'
$ws.Range("E17").Value = '  -0.70%  '
$ws.Range("D18").Value = '0.739'
$ws.Range("E18").Value = '  +1.54%  '
$ws.Range("D19").Value = '40.204.81'
$ws.Range("E19").Value = '  +0.57%  '
$ws.Range("E20").Value = '  +0.62%  '
$ws.Range("D21").Value = '11.40'
$ws.Range("E21").Value = '  +0.81%  '
$ws.Range("D22").Value = '5.80'
$ws.Range("E22").Value = '  -0.15%  '
$ws.Range("D23").Value = '65.86'
$ws.Range("E23").Value = '  +0.54%  '
$ws.Range("D24").Value = '236.38'
$ws.Range("E24").Value = '  +0.14%  '
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").Value = '2.52'
$ws.Range("E25").Value = '  +2.19%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").Value = '1.83'
$ws.Range("E27").Value = '  +1.01%  '
$ws.Range("D28").Value = '23.42'
$ws.Range("E28").Value = '  +3.07%  '
$ws.Range("E29").Value = '  +1.96%  '
$ws.Range("E30").Value = '  -1.04%  '
$ws.Range("D31").Value = '156.85'
$ws.Range("E31").Value = '  +0.40%  '
$ws.Range("D32").Value = '32.52'
$ws.Range("E32").Value = '  +2.57%  '
$ws.Range("E33").Value = '  +0.11%  '
$ws.Range("D34").Value = '5.00'
$ws.Range("E34").Value = '  +1.46%  '
$ws.Range("D35").Value = '0.0719'
$ws.Range("E35").Value = '  +0.43%  '
$ws.Range("E36").Value = '  +2.51%  '
$ws.Range("E37").Value = '  -0.58%  '
$ws.Range("E38").Value = '  +1.66%  '
$ws.Range("D39").Value = '0.103'
$ws.Range("E39").Value = '  +5.70%  '
$ws.Range("D40").Value = '1.74'
$ws.Range("E40").Value = '  +3.49%  '
$ws.Range("D41").Value = '15.64'
$ws.Range("E41").Value = '  -0.54%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").Value = '2.067.81'
$ws.Range("E43").Value = '  -2.65%  '
$ws.Range("D44").Value = '19.36'
$ws.Range("E44").Value = '  +7.37%  '
$ws.Range("D45").Value = '0.0272'
$ws.Range("E45").Value = '  +1.79%  '
$ws.Range("D46").Value = '10.07'
$ws.Range("E46").Value = '  +2.96%  '
$ws.Range("D47").Value = '2.84'
$ws.Range("E47").Value = '  +7.02%  '
$ws.Range("E48").Value = '  -12.86%  '
$ws.Range("D49").Value = '2.431.03'
$ws.Range("E49").Value = '  -0.10%  '
$ws.Range("D51").Value = '1.48'
$ws.Range("E51").Value = '  +1.11%  '

# Restore default (Normal) style on column D so no stray number formatting remains
$ws.Range("D2:D51").Style = "Normal"

Write-Output "Applied cryptos list update"
